$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a value while forcing text storage (avoids numeric auto-coercion
# for ambiguous strings like "521.18" or "1.00"), then reset the cell style so no
# new/residual number-format style is left applied to the cell.
function Set-TextValue {
    param($Cell, $Text)
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "58.258.86"
$ws.Range("E2").Value = "  -4.22%  "
Set-TextValue $ws.Range("D3") "2.643.58"
$ws.Range("E3").Value = "  -2.58%  "
$ws.Range("E4").Value = "  +0.08%  "
Set-TextValue $ws.Range("D5") "521.18"
$ws.Range("E5").Value = "  -1.02%  "
Set-TextValue $ws.Range("D6") "144.06"
$ws.Range("E6").Value = "  -0.53%  "
Set-TextValue $ws.Range("D7") "1.00"
$ws.Range("E7").Value = "  +0.29%  "
Set-TextValue $ws.Range("D8") "0.570"
$ws.Range("E8").Value = "  -1.58%  "
Set-TextValue $ws.Range("D9") "6.67"
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("E10").Value = "  -3.34%  "
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("E12").Value = "  +1.57%  "
Set-TextValue $ws.Range("D13") "3.108.22"
$ws.Range("E13").Value = "  -2.19%  "
Set-TextValue $ws.Range("D14") "58.278.02"
$ws.Range("E14").Value = "  -4.14%  "
Set-TextValue $ws.Range("D15") "20.86"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("E16").Value = "  -1.61%  "
Set-TextValue $ws.Range("D17") "2.646.66"
$ws.Range("E17").Value = "  -2.74%  "
Set-TextValue $ws.Range("D18") "337.94"
$ws.Range("E18").Value = "  -3.00%  "
$ws.Range("E19").Value = "  -2.67%  "
Set-TextValue $ws.Range("D20") "10.43"
$ws.Range("E20").Value = "  -1.44%  "
$ws.Range("E21").Value = "  -0.88%  "
Set-TextValue $ws.Range("D22") "1.00"
$ws.Range("E22").Value = "  +0.09%  "
Set-TextValue $ws.Range("D23") "64.58"
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  +0.59%  "
Set-TextValue $ws.Range("D25") "0.168"
$ws.Range("E25").Value = "  -2.08%  "
Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  +0.55%  "
Set-TextValue $ws.Range("D27") "0.0₃0795"
$ws.Range("E27").Value = "  -3.03%  "
$ws.Range("E28").Value = "  -2.94%  "
Set-TextValue $ws.Range("D29") "6.65"
$ws.Range("E29").Value = "  -2.96%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("E31").Value = "  -1.02%  "
Set-TextValue $ws.Range("D32") "152.70"
$ws.Range("E32").Value = "  +1.82%  "
Set-TextValue $ws.Range("D33") "18.85"
$ws.Range("E33").Value = "  -1.54%  "
$ws.Range("E34").Value = "  -3.19%  "
$ws.Range("B35").Value = "SuiNetwork"
$ws.Range("C35").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue $ws.Range("D35") "0.912"
$ws.Range("E35").Value = "  -3.64%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "1.18"
$ws.Range("E36").Value = "  -5.05%  "
Set-TextValue $ws.Range("D37") "0.858"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("E38").Value = "  -0.46%  "
$ws.Range("E39").Value = "  -5.34%  "
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("E42").Value = "  -0.80%  "
Set-TextValue $ws.Range("D43") "0.0969"
$ws.Range("E43").Value = "  -2.57%  "
Set-TextValue $ws.Range("D44") "269.59"
$ws.Range("E44").Value = "  -6.08%  "
Set-TextValue $ws.Range("D45") "19.40"
$ws.Range("E45").Value = "  -3.19%  "
Set-TextValue $ws.Range("D46") "0.0539"
$ws.Range("E46").Value = "  -0.42%  "
$ws.Range("E47").Value = "  +1.48%  "
Set-TextValue $ws.Range("D48") "2.035.15"
$ws.Range("E48").Value = "  -4.95%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D49") "0.0228"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D50") "4.67"
$ws.Range("E50").Value = "  -4.31%  "
Set-TextValue $ws.Range("D51") "18.29"
$ws.Range("E51").Value = "  -4.57%  "
